$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto market data.
# Some Price values are plain decimals with trailing zeros (e.g. "1.010", "90.00")
# that Excel would otherwise auto-convert to numbers and strip; force those cells
# to Text first so the literal formatted string is preserved, matching the source feed.

$ws.Range("D2").Value = "28.435.32"
$ws.Range("E2").Value = "  -3.55%  "
$ws.Range("D3").Value = "1.956.28"
$ws.Range("E3").Value = "  -1.85%  "
$ws.Range("E4").Value = "  -0.62%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "321.17"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.007"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.62%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4754"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -5.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4057"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -4.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.36"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.94%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08416"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -5.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.058"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -4.75%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.28"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -3.74%  "
$ws.Range("D13").Value = "1.950.04"
$ws.Range("E13").Value = "  -2.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.606"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -4.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.145"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -4.74%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.010"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "90.00"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -4.45%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001068"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -3.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06608"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.57"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -4.10%  "
$ws.Range("E21").Value = "  -0.58%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.818"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.88%  "
$ws.Range("D23").Value = "28.466.12"
$ws.Range("E23").Value = "  -3.54%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.52"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -4.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.292"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.36%  "
$ws.Range("D26").Value = "2.186.28"
$ws.Range("E26").Value = "  -2.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "155.16"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.17"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.90%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.914"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -5.55%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.155"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -6.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "123.59"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -3.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9784"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -7.53%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09587"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.444"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -6.49%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.597"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -3.91%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.661"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -3.47%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02339"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -4.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.905"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.17%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06212"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.51%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.250"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.72%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6214"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -4.82%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.15"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.82%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.007"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.58%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1920"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -5.70%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.356"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +3.59%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5953"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -5.96%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.99"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -4.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.061"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -6.58%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.394"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -3.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00000000328"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.98%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06828"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.79%  "
